$d = $word.ActiveDocument

# The logbook table ("Logboek") entries live in the 3rd w:tbl of the
# document body. Row 2 already holds the "20/11" entry ending in
# ". 3uur"; row 3 is the next (previously empty) entry that needs the
# new "22/11" date and its description.
$t = $d.Tables.Item(3)
$row = $t.Rows.Item(3)

# --- Date cell ("22/11") ---------------------------------------------
$dateCell = $row.Cells.Item(1)
$dateCell.Range.Text = "22/11"
# Re-fetch the cell range, trim the trailing cell-end mark, and apply
# the same font used throughout the logbook (Century Gothic) so the
# new run gets a matching <w:rPr><w:rFonts .../></w:rPr>.
$dateRange = $dateCell.Range
$dateTextRange = $d.Range($dateRange.Start, $dateRange.End - 1)
$dateTextRange.Font.Name = "Century Gothic"

# --- Description cell ---------------------------------------------------
$descCell = $row.Cells.Item(3)
$descText = "Maken + testen van het schema voor de lichtsturing, conclusie: Pin 3 van de 555 timer doet niks."
$descCell.Range.Text = $descText
$descRange = $descCell.Range
$descTextRange = $d.Range($descRange.Start, $descRange.End - 1)
$descTextRange.Font.Name = "Century Gothic"
